{"js": "// Update the date line and every \"A\u00d7B=\" practice problem in the table.\n// Each old value is unique in the document, so a scoped search-and-replace\n// (exact, case-sensitive, no wildcards) for each pair is unambiguous.\nconst replacements = [\n  [\"2025-11-17 Monday\", \"2025-11-18 Tuesday\"],\n  [\"788\u00d74=\", \"842\u00d74=\"],\n  [\"340\u00d76=\", \"290\u00d77=\"],\n  [\"746\u00d74=\", \"624\u00d76=\"],\n  [\"120\u00d73=\", \"744\u00d79=\"],\n  [\"785\u00d72=\", \"874\u00d78=\"],\n  [\"392\u00d72=\", \"501\u00d77=\"],\n  [\"180\u00d79=\", \"232\u00d78=\"],\n  [\"324\u00d72=\", \"403\u00d78=\"],\n  [\"257\u00d76=\", \"365\u00d73=\"],\n  [\"268\u00d74=\", \"802\u00d76=\"],\n  [\"895\u00d79=\", \"871\u00d79=\"],\n  [\"713\u00d79=\", \"737\u00d78=\"],\n  [\"910\u00d75=\", \"750\u00d77=\"],\n  [\"646\u00d77=\", \"411\u00d74=\"],\n  [\"325\u00d75=\", \"551\u00d76=\"],\n  [\"400\u00d74=\", \"896\u00d73=\"],\n  [\"903\u00d78=\", \"253\u00d77=\"],\n  [\"210\u00d79=\", \"637\u00d78=\"],\n  [\"345\u00d77=\", \"255\u00d77=\"],\n  [\"971\u00d78=\", \"207\u00d77=\"],\n  [\"339\u00d78=\", \"282\u00d76=\"],\n  [\"939\u00d74=\", \"885\u00d76=\"],\n  [\"553\u00d72=\", \"338\u00d74=\"],\n  [\"831\u00d76=\", \"470\u00d72=\"],\n  [\"312\u00d74=\", \"104\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every \"A\u00d7B=\" practice problem in the table.\n# Each old value is unique in the document, so a document-wide Find/Replace\n# (exact match, match case, no wildcards) for each pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-11-17 Monday\", \"2025-11-18 Tuesday\"),\n    @(\"788\u00d74=\", \"842\u00d74=\"),\n    @(\"340\u00d76=\", \"290\u00d77=\"),\n    @(\"746\u00d74=\", \"624\u00d76=\"),\n    @(\"120\u00d73=\", \"744\u00d79=\"),\n    @(\"785\u00d72=\", \"874\u00d78=\"),\n    @(\"392\u00d72=\", \"501\u00d77=\"),\n    @(\"180\u00d79=\", \"232\u00d78=\"),\n    @(\"324\u00d72=\", \"403\u00d78=\"),\n    @(\"257\u00d76=\", \"365\u00d73=\"),\n    @(\"268\u00d74=\", \"802\u00d76=\"),\n    @(\"895\u00d79=\", \"871\u00d79=\"),\n    @(\"713\u00d79=\", \"737\u00d78=\"),\n    @(\"910\u00d75=\", \"750\u00d77=\"),\n    @(\"646\u00d77=\", \"411\u00d74=\"),\n    @(\"325\u00d75=\", \"551\u00d76=\"),\n    @(\"400\u00d74=\", \"896\u00d73=\"),\n    @(\"903\u00d78=\", \"253\u00d77=\"),\n    @(\"210\u00d79=\", \"637\u00d78=\"),\n    @(\"345\u00d77=\", \"255\u00d77=\"),\n    @(\"971\u00d78=\", \"207\u00d77=\"),\n    @(\"339\u00d78=\", \"282\u00d76=\"),\n    @(\"939\u00d74=\", \"885\u00d76=\"),\n    @(\"553\u00d72=\", \"338\u00d74=\"),\n    @(\"831\u00d76=\", \"470\u00d72=\"),\n    @(\"312\u00d74=\", \"104\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
